$wb = $excel.ActiveWorkbook

# --- Sheet references ---
$wsInput = $wb.Worksheets.Item("Input")
$wsSummary = $wb.Worksheets.Item("Summary")
$wsRepay = $wb.Worksheets.Item("Repayment Schedule")
$wsTrans = $wb.Worksheets.Item("Transactions")

# --- Cell value / style updates ---

# Summary!F3 : new value, and number format reverts from "#,##0.00" (style 7)
# back to the plain wrap/vertical-centered style (style 3, like C3/D3 in the
# same row). Copy the format from D3 (already style 3) then overwrite the value.
$wsSummary.Range("D3").Copy() | Out-Null
$wsSummary.Range("F3").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$wsSummary.Range("F3").Value = 968.39

# Repayment Schedule!P2 -> O2 : the empty styled cell moves one column left.
$wsRepay.Range("P2").Copy() | Out-Null
$wsRepay.Range("O2").PasteSpecial(-4122) | Out-Null     # xlPasteFormats
$wsRepay.Range("P2").Clear() | Out-Null

# Transactions!A2 and A3 values
$wsTrans.Range("A2").Value = 229
$wsTrans.Range("A3").Value = 227

$excel.CutCopyMode = $false

# --- Selections / active sheet / active tab ---
# Final active sheet must be "Repayment Schedule" (tabSelected moves there,
# workbook activeTab index 2). Touch every sheet's selection along the way so
# each sheetView's <selection> reflects the new activeCell/sqref.

$wsInput.Activate()
$wsInput.Range("A2").Select() | Out-Null

$wsSummary.Activate()
$wsSummary.Range("C30").Select() | Out-Null

$wsTrans.Activate()
$wsTrans.Range("D2").Select() | Out-Null

$wsRepay.Activate()
$wsRepay.Range("G5").Select() | Out-Null
